$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117, shifting rows 117:141 down to 118:142
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new data record
$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(117, 3).Value = "Ñuble"
$ws.Cells.Item(117, 4).Value = 44504
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100112006
$ws.Cells.Item(117, 7).Value = "Repollo"
$ws.Cells.Item(117, 8).Value = "Crespo record"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 360
$ws.Cells.Item(117, 11).Value = 600
$ws.Cells.Item(117, 12).Value = 700
$ws.Cells.Item(117, 13).Value = 650
$ws.Cells.Item(117, 14).Value = "$/unidad"
$ws.Cells.Item(117, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(117, 16).Value = 650
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = "Hortaliza"
